$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 712, shifting all the
# following rows (712-782) down by one (to 713-783). This mirrors the
# weekly refresh where a new "Perejil" price record is prepended to the
# block and the rest of the previously recorded weeks shift down.
$ws.Rows.Item(712).Insert()

# Populate the newly inserted row with the latest week's record.
$ws.Range("A712").Value = 6
$ws.Range("B712").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C712").Value = "Metropolitana"
$ws.Range("D712").Value = 45132
$ws.Range("E712").Value = 13
$ws.Range("F712").Value = 100112044
$ws.Range("G712").Value = "Perejil"
$ws.Range("H712").Value = "Sin especificar"
$ws.Range("I712").Value = "Primera"
$ws.Range("J712").Value = 270
$ws.Range("K712").Value = 12000
$ws.Range("L712").Value = 13000
$ws.Range("M712").Value = 12444
$ws.Range("N712").Value = "$/docena de atados"
$ws.Range("O712").Value = "Región Metropolitana"
$ws.Range("P712").Value = 4148
$ws.Range("Q712").Value = 3
$ws.Range("R712").Value = "Hortaliza"
